# Update "Pais" worksheet with refreshed COVID-19 country/provincia stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 07:52"

# --- Israel (row 27): refreshed totals ---
$ws.Range("B27").Value = 291828
$ws.Range("C27").Value = 1335
$ws.Range("D27").Value = 238145
$ws.Range("E27").Value = 51700
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 1983

# --- Uzbekistan (row 59): refreshed totals ---
$ws.Range("B59").Value = 61205
$ws.Range("C59").Value = 107
$ws.Range("D59").Value = 58069
$ws.Range("E59").Value = 2629
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 507

# --- Kirguistan (row 68): refreshed totals ---
$ws.Range("B68").Value = 49528
$ws.Range("C68").Value = 298
$ws.Range("D68").Value = 44522
$ws.Range("E68").Value = 3916
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 1090

# --- Trinidad y Tobago's updated totals push it above Bahamas, Guinea
#     Ecuatorial and Surinam in the (descending) sort order, so rows
#     129-132 shift down by one and row 129 becomes Trinidad y Tobago ---
$ws.Range("A129").Value = "Trinidad yTobago"
$ws.Range("B129").Value = 5101
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 3252
$ws.Range("E129").Value = 1759
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 90

$ws.Range("A130").Value = "Bahamas"
$ws.Range("B130").Value = 5078
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 2900
$ws.Range("E130").Value = 2071
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 107

$ws.Range("A131").Value = "Guinea Ecuatorial"
$ws.Range("B131").Value = 5063
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 4894
$ws.Range("E131").Value = 86
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 83

$ws.Range("A132").Value = "Surinam"
$ws.Range("B132").Value = 5051
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 4845
$ws.Range("E132").Value = 99
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 107
